$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:F1) - matches authoring order of shared strings
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Fire rate"
$ws.Range("D1").Value = "Bullet speed"
$ws.Range("E1").Value = "Base ammo"
$ws.Range("F1").Value = "Reload speed"

# Row 2 - Cog / Pistol
$ws.Range("A2").Value = "Cog"
$ws.Range("B2").Value = "Pistol"
$ws.Range("C2").Value = 25
$ws.Range("D2").Value = 45
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 1

# G1 entered after row 2 (preserves shared-string authoring order)
$ws.Range("G1").Value = "Bullet per shot"

# Row 3 - Deek / Shotgun
$ws.Range("A3").Value = "Deek"
$ws.Range("B3").Value = "Shotgun"
$ws.Range("C3").Value = 15
$ws.Range("D3").Value = 50
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 3

# Header formatting: bold, 12pt, row height 15.75
$headerRange = $ws.Range("A1:G1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 12
$ws.Rows.Item(1).RowHeight = 15.75

# Column widths (B..H) to match the authored sheet's layout
$ws.Columns.Item(2).ColumnWidth = 8.666666666666666
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws.Columns.Item(4).ColumnWidth = 15.666666666666666
$ws.Columns.Item(5).ColumnWidth = 11.333333333333332
$ws.Columns.Item(6).ColumnWidth = 13.333333333333332
$ws.Columns.Item(7).ColumnWidth = 14.5
$ws.Columns.Item(8).ColumnWidth = 14.666666666666666

# Selection
$ws.Range("E9").Select()

# Page setup
$ws.PageSetup.Orientation = 1
